# New PO forecast model
# Updates "Weekly Quantity", "Monthly Trend" and "PO Forecast" sheets with
# new/refreshed forecast data.

$wb = $excel.ActiveWorkbook

# Date serial values are stored with the same "almost midnight" fractional
# offset used throughout the workbook (xx.99999999999).
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append 3 new weekly rows (16-18)
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyNewRows = @(
    @(16, 45662.99999999999, 17),
    @(17, 45669.99999999999, 25),
    @(18, 45676.99999999999, 12)
)

foreach ($row in $weeklyNewRows) {
    $r = $row[0]
    $wsWeekly.Cells.Item($r, 1).Value2 = $row[1]
    $wsWeekly.Cells.Item($r, 1).NumberFormat = $dateFormat
    $wsWeekly.Cells.Item($r, 2).Value2 = $row[2]
}

# ---------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append 1 new monthly row (7)
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Cells.Item(7, 1).Value2 = 45688.99999999999
$wsMonthly.Cells.Item(7, 1).NumberFormat = $dateFormat
$wsMonthly.Cells.Item(7, 2).Value2 = 54

# ---------------------------------------------------------------------
# Sheet 3: "PO Forecast" - refresh the whole forecast curve; rows 2-15
# keep their existing dates but get new forecast quantities, rows 16-23
# shift to later dates with new quantities, and rows 24-26 are brand new.
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$forecastRows = @(
    @(2,  45515.99999999999, 60),
    @(3,  45550.99999999999, 52),
    @(4,  45557.99999999999, 51),
    @(5,  45564.99999999999, 49),
    @(6,  45571.99999999999, 48),
    @(7,  45578.99999999999, 46),
    @(8,  45585.99999999999, 45),
    @(9,  45592.99999999999, 43),
    @(10, 45599.99999999999, 42),
    @(11, 45606.99999999999, 40),
    @(12, 45613.99999999999, 39),
    @(13, 45627.99999999999, 36),
    @(14, 45634.99999999999, 35),
    @(15, 45641.99999999999, 33),
    @(16, 45662.99999999999, 29),
    @(17, 45669.99999999999, 27),
    @(18, 45676.99999999999, 26),
    @(19, 45683.99999999999, 24),
    @(20, 45690.99999999999, 23),
    @(21, 45697.99999999999, 21),
    @(22, 45704.99999999999, 20),
    @(23, 45711.99999999999, 18),
    @(24, 45718.99999999999, 17),
    @(25, 45725.99999999999, 15),
    @(26, 45732.99999999999, 14)
)

foreach ($row in $forecastRows) {
    $r = $row[0]
    $wsForecast.Cells.Item($r, 1).Value2 = $row[1]
    $wsForecast.Cells.Item($r, 1).NumberFormat = $dateFormat
    $wsForecast.Cells.Item($r, 2).Value2 = $row[2]
}
